$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ideal-format")

# Replace the hard-coded population totals (N row) with formulas that sum
# the individual compartment values below (rows 16-22), for both the
# "o-o" (B) and "y-o"/"o-y" (C) scenario columns.
$ws.Range("B3").Formula = "=B16+B17+B18+B19+B20+B21+B22"
$ws.Range("C3").Formula = "=C16+C17+C18+C19+C20+C21+C22"

# Restore the cell selection/cursor position as last saved.
$ws.Range("C4").Select() | Out-Null
